# Applies the "New modelling cycle python code" edit:
#  - adds a new "New safety time:" column (G) to the "Times per aircraft"
#    sheet, computed as the general separation/safety time plus the
#    per-aircraft "Added safety time" (column F)
#  - widens the new columns to fit their content
#  - updates the remembered selections on both sheets

$wb = $excel.ActiveWorkbook

$wsGeneral = $wb.Worksheets.Item("General information")
$wsTimes   = $wb.Worksheets.Item("Times per aircraft")

# --- "Times per aircraft": new column G -------------------------------
$wsTimes.Range("G1").Value = "New safety time:"

for ($row = 2; $row -le 11; $row++) {
    $wsTimes.Range("G$row").Formula = "='General information'!`$B`$2+'Times per aircraft'!F$row"
}

# Match the widened columns (F/G) from the committed workbook as closely
# as this host's column-width quantization allows.
$wsTimes.Columns.Item(6).ColumnWidth = 14.92
$wsTimes.Columns.Item(7).ColumnWidth = 14.1

# --- Restore the selections recorded in the saved workbook ------------
$wsGeneral.Activate()
$wsGeneral.Range("B2").Select() | Out-Null

$wsTimes.Activate()
$wsTimes.Range("K8").Select() | Out-Null
